$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 371 (existing rows 371..395 shift down to 372..396,
# dimension grows from A1:R395 to A1:R396).
$ws.Rows("371").Insert()

# Populate the new row 371 with the new price-record observation.
$ws.Cells.Item(371, 1).Value = 4
$ws.Cells.Item(371, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(371, 3).Value = "Los Lagos"
$ws.Cells.Item(371, 4).Value = 44714
$ws.Cells.Item(371, 5).Value = 10
$ws.Cells.Item(371, 6).Value = 100114001
$ws.Cells.Item(371, 7).Value = "Papa"
$ws.Cells.Item(371, 8).Value = "Patagonia"
$ws.Cells.Item(371, 9).Value = "1a (guarda)"
$ws.Cells.Item(371, 10).Value = 300
$ws.Cells.Item(371, 11).Value = 7000
$ws.Cells.Item(371, 12).Value = 7500
$ws.Cells.Item(371, 13).Value = 7250
$ws.Cells.Item(371, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(371, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(371, 16).Value = 290
$ws.Cells.Item(371, 17).Value = 25
$ws.Cells.Item(371, 18).Value = "Hortaliza"
